$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# --- choices sheet: insert a new "display_name" column before the
#     existing label::language column (old column C -> becomes D) ---
$choices.Columns.Item(3).Insert()
$choices.Cells.Item(1, 3).Value = "display_name"
$choices.Cells.Item(1, 3).Font.Bold = $true

# Recompute the frozen pane split now that a column was inserted
# (freeze boundary shifts from C/D to D/E).
$choices.Activate()
[void]($excel.ActiveWindow.FreezePanes = $false)
[void]$choices.Range("E2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
[void]$choices.Range("C1").Select()

# --- zoom level: both survey and choices sheets are now shown at 150% ---
$survey.Activate()
$excel.ActiveWindow.Zoom = 150

$choices.Activate()
$excel.ActiveWindow.Zoom = 150

$survey.Activate()
